# Generate Report for Handoff
# Refresh the localization-status report: file "c45e7f02-d891-420d-8242-e56ae17f0ce6.md"
# has a new handback completed, so its handoff/handback timestamps are updated on the
# Overview sheet as well as the per-locale (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest Handoff Date" column for c45e7f02 row (row 6) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D6").Value = "2016-32-21 02:32:19"

# --- zh-cn sheet: "Latest Handoff Datetime" column for c45e7f02 row (row 6) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E6").Value = "2016-03-21 02:32:15"

# --- de-de sheet: "Latest Handoff Datetime" column for c45e7f02 row (row 6) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E6").Value = "2016-03-21 02:32:19"
